$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values for rows 2-25 (bus voltage magnitudes, per-unit),
# columns B-F and I-N (slack bus voltage changed from 1.05 to 1.02 pu,
# case re-run with 380 kV -> values recomputed accordingly).
$newValues = @{
    2 = @(1.02,1.014410119015922,1.021025038497922,1.016204316641466,1.030948338201141,1.027450432843899,1.019641515837839,1.023863588861564,1.019057185547152,1.033757872672202,1.010671207227263)
    3 = @(1.02,1.015207564669389,1.021596797218632,1.016876292130557,1.033037426978079,1.027601853828522,1.020074329994092,1.024242355375343,1.019534859168708,1.035651990433393,1.010815673652995)
    4 = @(1.02,1.015724132009217,1.021966947392641,1.017311995721368,1.034383820087226,1.027698209353829,1.020354314315732,1.024486926049972,1.019844157650965,1.036871878821053,1.010909087338569)
    5 = @(1.02,1.015941431482071,1.022122601710928,1.017495377361008,1.034948583317503,1.027738329410394,1.020472001011011,1.02458961953598,1.019974236549201,1.037383374693724,1.010948342444402)
    6 = @(1.02,1.015977924878406,1.022148739262417,1.017526180328978,1.035043336478438,1.027745043031227,1.020491760008001,1.024606854940642,1.019996080264428,1.037469178999971,1.010954932602889)
    7 = @(1.02,1.015727035048066,1.021969027084716,1.017314445245098,1.034391371398236,1.02769874696263,1.020355886925507,1.024488298733174,1.019845895576938,1.036878718706678,1.010909611930456)
    8 = @(1.02,1.014679501514761,1.021218228404856,1.016431228714467,1.031655488662263,1.027501943125076,1.019787802618121,1.023991701533587,1.019218573230526,1.034399202035357,1.010720043787694)
    9 = @(1.02,1.012838001654482,1.019896672377978,1.014881767232265,1.026791931409563,1.027142666878613,1.018786216972606,1.023112691926727,1.018114810213214,1.029984842608633,1.010385508447458)
    10 = @(1.02,1.011613348243571,1.019016656302753,1.013853495659026,1.023519062421532,1.026894691458589,1.018118161656682,1.022524052391527,1.017380133154249,1.027009854734637,1.01016216908785)
    11 = @(1.02,1.011083787162032,1.018635851297266,1.013409374171168,1.022094220216639,1.026785294036209,1.017828815040235,1.022268543276516,1.017062296202585,1.025713652097465,1.010065388782473)
    12 = @(1.02,1.010887194026731,1.018494441280972,1.013244578246248,1.021563781804672,1.02674435372838,1.017721328210787,1.022173542168589,1.016944280704375,1.02523094686992,1.010029429543413)
    13 = @(1.02,1.010929358954317,1.018524772482747,1.013279919792719,1.021677617022722,1.026753149398519,1.017744384972276,1.022193924464245,1.016969593457847,1.025334545346605,1.010037143403989)
    14 = @(1.02,1.011067534474274,1.01862416152377,1.013395748601757,1.022050398472011,1.026781916132431,1.017819930358799,1.022260692364869,1.017052540117036,1.025673776983174,1.010062416598043)
    15 = @(1.02,1.011152683501612,1.018685403437026,1.013467137193669,1.022279922825819,1.026799599778013,1.01786647497441,1.022301817867132,1.017103652051446,1.025882623736908,1.010077986832592)
    16 = @(1.02,1.011648508761574,1.019041934356778,1.013882994390362,1.023613459842174,1.026901909077366,1.018137363120905,1.022540996549757,1.017401232975864,1.027095707719785,1.01016859056173)
    17 = @(1.02,1.011959720814994,1.019265643614679,1.014144153223269,1.024447873781138,1.026965542567823,1.01830726460604,1.022690859984938,1.017587973959977,1.027854474089862,1.01022540453422)
    18 = @(1.02,1.012141315132708,1.019396153301059,1.014296591397033,1.024933835927349,1.027002463866998,1.018406358068182,1.022778212546049,1.017696924077793,1.028296278835597,1.010258536115762)
    19 = @(1.02,1.01220324587628,1.01944065775687,1.014348587263929,1.025099412482374,1.027015020035162,1.018440145139796,1.022807987288206,1.017734077852539,1.02844679310756,1.01026983192647)
    20 = @(1.02,1.01192632350185,1.019241639255446,1.014116122105359,1.02435842559609,1.026958735480329,1.01828903652545,1.022674787287032,1.017567935594435,1.027773145615001,1.010219309659986)
    21 = @(1.02,1.011026842199099,1.018594892876223,1.013361635185783,1.021940656733222,1.026773453483559,1.017797684414071,1.022241033466281,1.017028113191966,1.025573916167745,1.010054974566571)
    22 = @(1.02,1.010461936195346,1.018188477464799,1.01288824679529,1.020413614220084,1.026655192634252,1.017488690694026,1.021967773473415,1.016688956297232,1.024183994236191,1.009951588630154)
    23 = @(1.02,1.010761343179222,1.018403904913134,1.013139104916922,1.021223794288868,1.026718052883385,1.017652499729762,1.022112685073276,1.016868725750584,1.024921510134177,1.010006401268262)
    24 = @(1.02,1.011941414095561,1.019252485718773,1.014128787825695,1.024398845657799,1.026961811912871,1.01829727303299,1.02268205003378,1.017576989985371,1.027809896832652,1.01022206369318)
    25 = @(1.02,1.013313546977621,1.020238149729048,1.015281516912618,1.028054514175349,1.027237034275327,1.019045211983765,1.023340402269977,1.018399957652103,1.031131585945647,1.010472050595858)
}

$cols = @("B","C","D","E","F","I","J","K","L","M","N")
foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}
